$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2008 and 2009 rows; every later row shifts up by two.
$ws.Rows("2:3").Delete()

# Append the new 2021 row at the bottom (now row 13), matching the
# format of the row above it.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 391.5
$ws.Range("C13").Value = 821.8
$ws.Range("D13").Value = 10559.9
$ws.Range("F13").Value = 58979
$ws.Range("G13").Value = 93318.60000000001
$ws.Range("H13").Value = 28528.5
$ws.Range("I13").Value = 9824.200000000001
$ws.Range("J13").Value = 21947
$ws.Range("K13").Value = 36242.4
$ws.Range("L13").Value = 75850
$ws.Range("N13").Value = 1503.9
$ws.Range("O13").Value = 8558.1
$ws.Range("P13").Value = 14195.6
$ws.Range("Q13").Value = 112092.4
$ws.Range("R13").Value = 214.4
$ws.Range("S13").Value = 4085.2
